$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.754972666666667
$ws.Range("H2").Value = 5.264918
$ws.Range("I2").Value = 0.5110994274238188
$ws.Range("J2").Value = 0.5110994274238188
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.970048
$ws.Range("N2").Value = 26.910144
$ws.Range("O2").Value = 0.487108783009476
$ws.Range("P2").Value = 0.4871087830094759
$ws.Range("Q2").Value = 15.742189058688
$ws.Range("R2").Value = 141.679701528192
$ws.Range("S2").Value = 0.2489610200892563
$ws.Range("T2").Value = 0.2489610200892563
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.754972666666667
$ws.Range("H3").Value = 5.264918
$ws.Range("I3").Value = 0.5110994274238188
$ws.Range("J3").Value = 0.5110994274238188
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.012070666666666
$ws.Range("N3").Value = 27.036212
$ws.Range("O3").Value = 0.489390778604016
$ws.Range("P3").Value = 0.489390778604016
$ws.Range("Q3").Value = 15.81593769006844
$ws.Range("R3").Value = 142.343439210616
$ws.Range("S3").Value = 0.2501273467310094
$ws.Range("T3").Value = 0.2501273467310094
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.754972666666667
$ws.Range("H4").Value = 5.264918
$ws.Range("I4").Value = 0.5110994274238188
$ws.Range("J4").Value = 0.5110994274238188
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4327576666666667
$ws.Range("N4").Value = 1.298273
$ws.Range("O4").Value = 0.02350043838650813
$ws.Range("P4").Value = 0.02350043838650813
$ws.Range("Q4").Value = 0.7594778762904444
$ws.Range("R4").Value = 6.835300886614
$ws.Range("S4").Value = 0.01201106060355303
$ws.Range("T4").Value = 0.01201106060355303
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.678748
$ws.Range("H5").Value = 5.036244
$ws.Range("I5").Value = 0.4889005725761812
$ws.Range("J5").Value = 0.4889005725761812
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.970048
$ws.Range("N5").Value = 26.910144
$ws.Range("O5").Value = 0.487108783009476
$ws.Range("P5").Value = 0.4871087830094759
$ws.Range("Q5").Value = 15.058450139904
$ws.Range("R5").Value = 135.526051259136
$ws.Range("S5").Value = 0.2381477629202196
$ws.Range("T5").Value = 0.2381477629202196
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.678748
$ws.Range("H6").Value = 5.036244
$ws.Range("I6").Value = 0.4889005725761812
$ws.Range("J6").Value = 0.4889005725761812
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.012070666666666
$ws.Range("N6").Value = 27.036212
$ws.Range("O6").Value = 0.489390778604016
$ws.Range("P6").Value = 0.489390778604016
$ws.Range("Q6").Value = 15.12899560752533
$ws.Range("R6").Value = 136.160960467728
$ws.Range("S6").Value = 0.2392634318730065
$ws.Range("T6").Value = 0.2392634318730066
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.678748
$ws.Range("H7").Value = 5.036244
$ws.Range("I7").Value = 0.4889005725761812
$ws.Range("J7").Value = 0.4889005725761812
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.4327576666666667
$ws.Range("N7").Value = 1.298273
$ws.Range("O7").Value = 0.02350043838650813
$ws.Range("P7").Value = 0.02350043838650813
$ws.Range("Q7").Value = 0.7264910674013333
$ws.Range("R7").Value = 6.538419606612
$ws.Range("S7").Value = 0.01148937778295509
$ws.Range("T7").Value = 0.01148937778295509
